# Normalize the "Recorded By" (column G) values so that any "System"/"system"
# tokens are moved to the front of the comma-separated list, while preserving
# the relative order of all other tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $newParts = $systemParts + $otherParts
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
